# Add two new columns (AD: diceSetupHexOffset, AE: diceSetupClockwise) and
# two new data rows (6 and 7) to Sheet1, matching a new "dice setup" feature.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New headers (row 1) -------------------------------------------------
# (NB: these two header strings carry a leading UTF-8 BOM (U+FEFF) in the
# source workbook's sharedStrings table, so we reproduce it verbatim.)
$ws.Cells.Item(1, 30).Value = "﻿diceSetupHexOffset"
$ws.Cells.Item(1, 31).Value = "﻿diceSetupClockwise"

# --- Fill AD/AE for existing rows 2-5 with defaults (0 / FALSE) --------
for ($r = 2; $r -le 5; $r++) {
    $ws.Cells.Item($r, 30).Value = 0
    $ws.Cells.Item($r, 31).Value = $false
}

# --- Tile/port resource assignments for the two new games --------------
$row6 = @("sheep","wheat","brick","brick","brick","rock","rock","sheep","wood","sheep","wood","wheat","rock","wheat","wood","wood","sheep","wheat","desert","brick","wood","desert","wheat","rock","desert","desert","sheep","desert")
$row7 = @("wood","brick","wood","rock","sheep","wheat","brick","wheat","rock","wood","sheep","rock","sheep","sheep","wheat","desert","brick","wood","wheat","brick","rock","sheep","desert","desert","desert","wood","desert","wheat")

# Row 6 -> Game 5
$ws.Cells.Item(6, 1).Value = 5
for ($i = 0; $i -lt $row6.Length; $i++) {
    $ws.Cells.Item(6, 2 + $i).Value = $row6[$i]
}
$ws.Cells.Item(6, 30).Value = 0
$ws.Cells.Item(6, 31).Value = $false

# Row 7 -> Game 6
$ws.Cells.Item(7, 1).Value = 6
for ($i = 0; $i -lt $row7.Length; $i++) {
    $ws.Cells.Item(7, 2 + $i).Value = $row7[$i]
}
$ws.Cells.Item(7, 30).Value = -2
$ws.Cells.Item(7, 31).Value = $false

# --- Column widths for the new columns (AD/AE), matching Excel's autofit
# (best-fit widths as Excel would compute them for these header strings)
$ws.Columns.Item(30).ColumnWidth = 16.666666666666668
$ws.Columns.Item(31).ColumnWidth = 16.498697916666668

# --- View state: scroll right so the new columns are visible, then select
# the last-entered cell (AD7), matching where the author's cursor ended up.
$ws.Application.ActiveWindow.ScrollColumn = 19
$ws.Range("AD7").Select()

$wb.Save()
